$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.489.42"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "2.638.61"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "326.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("E14").Value = "  +2.35%  "
$ws.Range("D15").Value = "3.051.84"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "2.648.15"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.853"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").Value = "49.488.73"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "268.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.24%  "
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0809"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.38%  "
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "128.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.65%  "
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0331"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.62%  "
$ws.Range("D45").Value = "2.057.74"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.66%  "
$ws.Range("E48").Value = "  -5.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.97%  "
